$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated stats for existing rows (row number -> values for columns D..S)
$updates = @{
    3  = @(29, 16, 14.06, 657, 321, 59, 71, 3671, 85, 133, 228, 2.931034482758621, 4.586206896551724, 7.862068965517241, 22.66, 29.09)
    4  = @(146, 116, 72.0156666666667, 2896, 1252, 399, 392, 29178, 1159, 886, 1100, 7.938356164383562, 6.068493150684931, 7.534246575342466, 19.84, 29.6)
    5  = @(50, 42, 26.28183333333333, 2866, 1391, 388, 363, 4617, 271, 203, 575, 5.42, 4.06, 11.5, 57.32, 31.54)
    8  = @(31, 14, 15.15316666666667, 2253, 957, 342, 219, 1158, 95, 222, 425, 3.064516129032258, 7.161290322580645, 13.70967741935484, 72.68000000000001, 29.33)
    10 = @(65, 7, 32.70866666666667, 5181, 2958, 633, 701, 1246, 154, 191, 1079, 2.369230769230769, 2.938461538461539, 16.6, 79.70999999999999, 30.19)
}

foreach ($rowNum in $updates.Keys) {
    $values = $updates[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column D is index 4
        $col = 4 + $i
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}

# New rows 12 and 13
$newRows = @{
    12 = "chguizou"
    13 = "zyradelevingne"
}

foreach ($rowNum in $newRows.Keys) {
    $ws.Cells.Item($rowNum, 1).Value = $newRows[$rowNum]
    for ($col = 2; $col -le 19; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = 0
    }
}
